# Generate Report for Handoff
# Adds two newly-discovered source files to the localization status report:
#   1ef2978a-9a86-4a10-9a02-5ad409f9842e.md
#   d0f8de88-3a56-4619-98c5-ea3770e9334b.md
# Both are freshly queued ("Ready for handoff") with generated handoff
# artifacts (.xlf) for the zh-cn and de-de locales, but have not yet been
# handed back (placeholder 0001-01-01 00:00:00 datetime, "Include" dependency).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# File / commit identifiers for the two new entries
# ---------------------------------------------------------------------------
$file1 = "1ef2978a-9a86-4a10-9a02-5ad409f9842e"
$file2 = "d0f8de88-3a56-4619-98c5-ea3770e9334b"

$xlf1 = "52a7e53c706c2d94ef67f4470b1d6e79b755288e"
$xlf2 = "8b2b0290bbf72e62524b6af51102604367572382"

$status = "Ready for handoff"
$noHandback = "0001-01-01 00:00:00"
$dependency = "Include"

$handoffDateOverview1 = "2016-30-20 00:30:32"
$handoffDateOverview2 = "2016-30-20 00:30:32"

$handoffDateZh1 = "2016-03-20 00:30:28"
$handoffDateZh2 = "2016-03-20 00:30:28"

$handoffDateDe1 = "2016-03-20 00:30:32"
$handoffDateDe2 = "2016-03-20 00:30:32"

# ---------------------------------------------------------------------------
# Overview sheet: one row per source file
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "$file1.md"
$ws.Range("B6").Value = $status
$ws.Range("C6").Value = $status
$ws.Range("D6").Value = $handoffDateOverview1

$ws.Range("A7").Value = "$file2.md"
$ws.Range("B7").Value = $status
$ws.Range("C7").Value = $status
$ws.Range("D7").Value = $handoffDateOverview2

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/13789de84663b629c086ee53f2bea8d44a7c8e78/e2e/$file1.md", "", "", "$file1.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0063f2b3261162024fee2f0f5eb15ee4112ad4b5/e2e/$file2.md", "", "", "$file2.md")

# ---------------------------------------------------------------------------
# zh-cn sheet: locale detail rows
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A6").Value = "$file1.md"
$wsZh.Range("B6").Value = ".md"
$wsZh.Range("C6").Value = $status
$wsZh.Range("D6").Value = "$file1.$xlf1.zh-cn.xlf"
$wsZh.Range("E6").Value = $handoffDateZh1
$wsZh.Range("H6").Value = $noHandback
$wsZh.Range("I6").Value = $dependency

$wsZh.Range("A7").Value = "$file2.md"
$wsZh.Range("B7").Value = ".md"
$wsZh.Range("C7").Value = $status
$wsZh.Range("D7").Value = "$file2.$xlf2.zh-cn.xlf"
$wsZh.Range("E7").Value = $handoffDateZh2
$wsZh.Range("H7").Value = $noHandback
$wsZh.Range("I7").Value = $dependency

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/13789de84663b629c086ee53f2bea8d44a7c8e78/e2e/$file1.md", "", "", "$file1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d237e55a685739aeed9a209ca693bc0f6a93a9ac/e2e/$file1.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f17c7af8d5f69fda16179f614357a49351d044a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file1.$xlf1.zh-cn.xlf", "", "", "$file1.$xlf1.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0063f2b3261162024fee2f0f5eb15ee4112ad4b5/e2e/$file2.md", "", "", "$file2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6db9b849be4a9e8c863572ea32f5c41e11c76ea8/e2e/$file2.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d340e598a958699bcb52a9389051c20d72ef6cb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file2.$xlf2.zh-cn.xlf", "", "", "$file2.$xlf2.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet: locale detail rows
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A6").Value = "$file1.md"
$wsDe.Range("B6").Value = ".md"
$wsDe.Range("C6").Value = $status
$wsDe.Range("D6").Value = "$file1.$xlf1.de-de.xlf"
$wsDe.Range("E6").Value = $handoffDateDe1
$wsDe.Range("H6").Value = $noHandback
$wsDe.Range("I6").Value = $dependency

$wsDe.Range("A7").Value = "$file2.md"
$wsDe.Range("B7").Value = ".md"
$wsDe.Range("C7").Value = $status
$wsDe.Range("D7").Value = "$file2.$xlf2.de-de.xlf"
$wsDe.Range("E7").Value = $handoffDateDe2
$wsDe.Range("H7").Value = $noHandback
$wsDe.Range("I7").Value = $dependency

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/13789de84663b629c086ee53f2bea8d44a7c8e78/e2e/$file1.md", "", "", "$file1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/523d3cb120ecf683633b6ba9264eb9eb2ae1d48b/e2e/$file1.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f17c7af8d5f69fda16179f614357a49351d044a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file1.$xlf1.de-de.xlf", "", "", "$file1.$xlf1.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0063f2b3261162024fee2f0f5eb15ee4112ad4b5/e2e/$file2.md", "", "", "$file2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4fd44c035ba1a9ce2a34363854786d4173df8986/e2e/$file2.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d340e598a958699bcb52a9389051c20d72ef6cb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file2.$xlf2.de-de.xlf", "", "", "$file2.$xlf2.de-de.xlf")

Write-Output "Report generated: added $file1 and $file2 across Overview, zh-cn, de-de sheets."
